$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.028117330279534
$ws.Cells.Item(2, 4).Value = 1.029094086837847
$ws.Cells.Item(2, 5).Value = 1.037278728939197
$ws.Cells.Item(2, 6).Value = 1.046308381459464
$ws.Cells.Item(2, 9).Value = 1.031046868218264
$ws.Cells.Item(2, 10).Value = 1.033271429457739
$ws.Cells.Item(2, 11).Value = 1.031908996946797
$ws.Cells.Item(2, 12).Value = 1.040070068345771
$ws.Cells.Item(2, 13).Value = 1.049074180890877
$ws.Cells.Item(2, 14).Value = 1.034738793405458
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.029335688008883
$ws.Cells.Item(3, 4).Value = 1.029461717418037
$ws.Cells.Item(3, 5).Value = 1.0383603432366
$ws.Cells.Item(3, 6).Value = 1.047474702360303
$ws.Cells.Item(3, 9).Value = 1.031085851222385
$ws.Cells.Item(3, 10).Value = 1.034128881572542
$ws.Cells.Item(3, 11).Value = 1.032085888028109
$ws.Cells.Item(3, 12).Value = 1.04096073765349
$ws.Cells.Item(3, 13).Value = 1.050051196407676
$ws.Cells.Item(3, 14).Value = 1.035597463200616
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.030124201309196
$ws.Cells.Item(4, 4).Value = 1.029699380272863
$ws.Cells.Item(4, 5).Value = 1.03906056858344
$ws.Cells.Item(4, 6).Value = 1.048229547640932
$ws.Cells.Item(4, 9).Value = 1.031109525754181
$ws.Cells.Item(4, 10).Value = 1.034683387065572
$ws.Cells.Item(4, 11).Value = 1.032199383715392
$ws.Cells.Item(4, 12).Value = 1.041536817977136
$ws.Cells.Item(4, 13).Value = 1.050682973457008
$ws.Cells.Item(4, 14).Value = 1.036152756155049
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.030455731153041
$ws.Cells.Item(5, 4).Value = 1.029799239420322
$ws.Cells.Item(5, 5).Value = 1.039355027619347
$ws.Cells.Item(5, 6).Value = 1.048546923167023
$ws.Cells.Item(5, 9).Value = 1.031119106756078
$ws.Cells.Item(5, 10).Value = 1.03491642471103
$ws.Cells.Item(5, 11).Value = 1.032246865114289
$ws.Cells.Item(5, 12).Value = 1.041778944939222
$ws.Cells.Item(5, 13).Value = 1.05094847328602
$ws.Cells.Item(5, 14).Value = 1.036386124740693
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.030511398806059
$ws.Cells.Item(6, 4).Value = 1.029816002968292
$ws.Cells.Item(6, 5).Value = 1.039404473555316
$ws.Cells.Item(6, 6).Value = 1.048600214200099
$ws.Cells.Item(6, 9).Value = 1.031120693627789
$ws.Cells.Item(6, 10).Value = 1.034955548300859
$ws.Cells.Item(6, 11).Value = 1.032254823798151
$ws.Cells.Item(6, 12).Value = 1.041819595790102
$ws.Cells.Item(6, 13).Value = 1.050993046075255
$ws.Cells.Item(6, 14).Value = 1.036425303890507
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.030128631072166
$ws.Cells.Item(7, 4).Value = 1.029700714811595
$ws.Cells.Item(7, 5).Value = 1.039064502827195
$ws.Cells.Item(7, 6).Value = 1.048233788275682
$ws.Cells.Item(7, 9).Value = 1.031109655237663
$ws.Cells.Item(7, 10).Value = 1.034686501225271
$ws.Cells.Item(7, 11).Value = 1.032200019078088
$ws.Cells.Item(7, 12).Value = 1.041540053513255
$ws.Cells.Item(7, 13).Value = 1.050686521467867
$ws.Cells.Item(7, 14).Value = 1.036155874737211
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.028529049151892
$ws.Cells.Item(8, 4).Value = 1.029218373427148
$ws.Cells.Item(8, 5).Value = 1.037644193893993
$ws.Cells.Item(8, 6).Value = 1.046702512918739
$ws.Cells.Item(8, 9).Value = 1.031060363253687
$ws.Cells.Item(8, 10).Value = 1.033561276940233
$ws.Cells.Item(8, 11).Value = 1.031968977111784
$ws.Cells.Item(8, 12).Value = 1.040371124852199
$ws.Cells.Item(8, 13).Value = 1.049404454669048
$ws.Cells.Item(8, 14).Value = 1.035029052504624
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.025711478988682
$ws.Cells.Item(9, 4).Value = 1.028366837238586
$ws.Cells.Item(9, 5).Value = 1.035144059201943
$ws.Cells.Item(9, 6).Value = 1.044005380268634
$ws.Cells.Item(9, 9).Value = 1.030961662574238
$ws.Cells.Item(9, 10).Value = 1.031575968642434
$ws.Cells.Item(9, 11).Value = 1.031554518101667
$ws.Cells.Item(9, 12).Value = 1.038309431699046
$ws.Cells.Item(9, 13).Value = 1.047142062722855
$ws.Cells.Item(9, 14).Value = 1.033040924841324
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.02383371633209
$ws.Cells.Item(10, 4).Value = 1.027798194150698
$ws.Cells.Item(10, 5).Value = 1.03347901531746
$ws.Cells.Item(10, 6).Value = 1.042208023366235
$ws.Cells.Item(10, 9).Value = 1.030887940287749
$ws.Cells.Item(10, 10).Value = 1.030250674100194
$ws.Cells.Item(10, 11).Value = 1.031273355053899
$ws.Cells.Item(10, 12).Value = 1.036933645319921
$ws.Cells.Item(10, 13).Value = 1.04563158705792
$ws.Cells.Item(10, 14).Value = 1.03171374822882
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.02302074340784
$ws.Cells.Item(11, 4).Value = 1.027551764079307
$ws.Cells.Item(11, 5).Value = 1.032758425765578
$ws.Cells.Item(11, 6).Value = 1.041429908268601
$ws.Cells.Item(11, 9).Value = 1.030854147116114
$ws.Cells.Item(11, 10).Value = 1.029676376118096
$ws.Cells.Item(11, 11).Value = 1.031150471825462
$ws.Cells.Item(11, 12).Value = 1.036337588121981
$ws.Cells.Item(11, 13).Value = 1.044976998286283
$ws.Cells.Item(11, 14).Value = 1.031138634677715
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.022718783252297
$ws.Cells.Item(12, 4).Value = 1.027460200341883
$ws.Cells.Item(12, 5).Value = 1.032490823383385
$ws.Cells.Item(12, 6).Value = 1.041140903307517
$ws.Cells.Item(12, 9).Value = 1.030841314388672
$ws.Cells.Item(12, 10).Value = 1.029462989345189
$ws.Cells.Item(12, 11).Value = 1.031104657827524
$ws.Cells.Item(12, 12).Value = 1.036116134809038
$ws.Cells.Item(12, 13).Value = 1.044733772266124
$ws.Cells.Item(12, 14).Value = 1.030924944871114
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.022783554128012
$ws.Cells.Item(13, 4).Value = 1.02747984233842
$ws.Cells.Item(13, 5).Value = 1.032548222458748
$ws.Cells.Item(13, 6).Value = 1.041202894893944
$ws.Cells.Item(13, 9).Value = 1.030844079732356
$ws.Cells.Item(13, 10).Value = 1.029508764582398
$ws.Cells.Item(13, 11).Value = 1.031114492744828
$ws.Cells.Item(13, 12).Value = 1.03616363962878
$ws.Cells.Item(13, 13).Value = 1.04478594882802
$ws.Cells.Item(13, 14).Value = 1.030970785114411
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.022995783011939
$ws.Cells.Item(14, 4).Value = 1.027544195964923
$ws.Cells.Item(14, 5).Value = 1.032736304520526
$ws.Cells.Item(14, 6).Value = 1.041406018594324
$ws.Cells.Item(14, 9).Value = 1.030853092075995
$ws.Cells.Item(14, 10).Value = 1.029658738873926
$ws.Cells.Item(14, 11).Value = 1.031146688280687
$ws.Cells.Item(14, 12).Value = 1.0363192837718
$ws.Cells.Item(14, 13).Value = 1.044956894838777
$ws.Cells.Item(14, 14).Value = 1.031120972386636
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.023126546038408
$ws.Cells.Item(15, 4).Value = 1.027583842622838
$ws.Cells.Item(15, 5).Value = 1.032852195557168
$ws.Cells.Item(15, 6).Value = 1.041531172639025
$ws.Cells.Item(15, 9).Value = 1.030858607736802
$ws.Cells.Item(15, 10).Value = 1.02975113407601
$ws.Cells.Item(15, 11).Value = 1.031166502563614
$ws.Cells.Item(15, 12).Value = 1.036415174455283
$ws.Cells.Item(15, 13).Value = 1.045062209340555
$ws.Cells.Item(15, 14).Value = 1.031213498800504
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.023887672791227
$ws.Cells.Item(16, 4).Value = 1.027814544778102
$ws.Cells.Item(16, 5).Value = 1.033526846487262
$ws.Cells.Item(16, 6).Value = 1.042259667404828
$ws.Cells.Item(16, 9).Value = 1.030890143666203
$ws.Cells.Item(16, 10).Value = 1.030288779052265
$ws.Cells.Item(16, 11).Value = 1.031281486543706
$ws.Cells.Item(16, 12).Value = 1.036973196566424
$ws.Cells.Item(16, 13).Value = 1.045675018432354
$ws.Cells.Item(16, 14).Value = 1.031751907294294
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.024365135405186
$ws.Cells.Item(17, 4).Value = 1.027959204975412
$ws.Cells.Item(17, 5).Value = 1.033950139505025
$ws.Cells.Item(17, 6).Value = 1.042716672818436
$ws.Cells.Item(17, 9).Value = 1.030909424877651
$ws.Cells.Item(17, 10).Value = 1.030625911691083
$ws.Cells.Item(17, 11).Value = 1.031353309221757
$ws.Cells.Item(17, 12).Value = 1.037323139225624
$ws.Cells.Item(17, 13).Value = 1.046059271311329
$ws.Cells.Item(17, 14).Value = 1.032089518700128
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.024643642116899
$ws.Cells.Item(18, 4).Value = 1.02804356301878
$ws.Cells.Item(18, 5).Value = 1.034197076435427
$ws.Cells.Item(18, 6).Value = 1.04298325112517
$ws.Cells.Item(18, 9).Value = 1.030920490720905
$ws.Cells.Item(18, 10).Value = 1.030822513410697
$ws.Cells.Item(18, 11).Value = 1.031395092313457
$ws.Cells.Item(18, 12).Value = 1.037527222844882
$ws.Cells.Item(18, 13).Value = 1.046283347210572
$ws.Cells.Item(18, 14).Value = 1.032286399616735
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.024738607621658
$ws.Cells.Item(19, 4).Value = 1.028072323511874
$ws.Cells.Item(19, 5).Value = 1.034281281906872
$ws.Cells.Item(19, 6).Value = 1.043074149983076
$ws.Cells.Item(19, 9).Value = 1.030924233241339
$ws.Cells.Item(19, 10).Value = 1.030889542424505
$ws.Cells.Item(19, 11).Value = 1.031409320595589
$ws.Cells.Item(19, 12).Value = 1.037596804689626
$ws.Cells.Item(19, 13).Value = 1.046359742461098
$ws.Cells.Item(19, 14).Value = 1.032353523819431
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.024313907082279
$ws.Cells.Item(20, 4).Value = 1.027943686339297
$ws.Cells.Item(20, 5).Value = 1.033904720324707
$ws.Cells.Item(20, 6).Value = 1.04266763894386
$ws.Cells.Item(20, 9).Value = 1.030907374857855
$ws.Cells.Item(20, 10).Value = 1.030589744908832
$ws.Cells.Item(20, 11).Value = 1.031345614682738
$ws.Cells.Item(20, 12).Value = 1.037285597010637
$ws.Cells.Item(20, 13).Value = 1.046018050032535
$ws.Cells.Item(20, 14).Value = 1.032053300556898
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.022933286523324
$ws.Cells.Item(21, 4).Value = 1.027525246200477
$ws.Cells.Item(21, 5).Value = 1.032680917487039
$ws.Cells.Item(21, 6).Value = 1.041346203140201
$ws.Cells.Item(21, 9).Value = 1.03085044590172
$ws.Cells.Item(21, 10).Value = 1.029614577041358
$ws.Cells.Item(21, 11).Value = 1.031137212172323
$ws.Cells.Item(21, 12).Value = 1.036273451871767
$ws.Cells.Item(21, 13).Value = 1.044906557788065
$ws.Cells.Item(21, 14).Value = 1.0310767478392
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.022065316446286
$ws.Cells.Item(22, 4).Value = 1.027261991570284
$ws.Cells.Item(22, 5).Value = 1.031911791277803
$ws.Cells.Item(22, 6).Value = 1.040515488535969
$ws.Cells.Item(22, 9).Value = 1.030813030192729
$ws.Cells.Item(22, 10).Value = 1.029001062395929
$ws.Cells.Item(22, 11).Value = 1.031005200045155
$ws.Cells.Item(22, 12).Value = 1.035636779652114
$ws.Cells.Item(22, 13).Value = 1.04420724005106
$ws.Cells.Item(22, 14).Value = 1.030462361932603
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.02252543695867
$ws.Cells.Item(23, 4).Value = 1.027401562781779
$ws.Cells.Item(23, 5).Value = 1.032319489035313
$ws.Cells.Item(23, 6).Value = 1.040955854679544
$ws.Cells.Item(23, 9).Value = 1.030833018519857
$ws.Cells.Item(23, 10).Value = 1.029326335328669
$ws.Cells.Item(23, 11).Value = 1.031075274765231
$ws.Cells.Item(23, 12).Value = 1.03597432012485
$ws.Cells.Item(23, 13).Value = 1.044578007317738
$ws.Cells.Item(23, 14).Value = 1.030788096790216
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.024337054922801
$ws.Cells.Item(24, 4).Value = 1.027950698604475
$ws.Cells.Item(24, 5).Value = 1.033925243182189
$ws.Cells.Item(24, 6).Value = 1.042689795195606
$ws.Cells.Item(24, 9).Value = 1.030908301731581
$ws.Cells.Item(24, 10).Value = 1.030606087250922
$ws.Cells.Item(24, 11).Value = 1.031349091853373
$ws.Cells.Item(24, 12).Value = 1.037302560820606
$ws.Cells.Item(24, 13).Value = 1.046036676315802
$ws.Cells.Item(24, 14).Value = 1.032069666106989
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.026439770752477
$ws.Cells.Item(25, 4).Value = 1.028587156053423
$ws.Cells.Item(25, 5).Value = 1.035790097948365
$ws.Cells.Item(25, 6).Value = 1.044702520712501
$ws.Cells.Item(25, 9).Value = 1.030988578257579
$ws.Cells.Item(25, 10).Value = 1.032089523338502
$ws.Cells.Item(25, 11).Value = 1.031662527136736
$ws.Cells.Item(25, 12).Value = 1.038842658351631
$ws.Cells.Item(25, 13).Value = 1.047727332251492
$ws.Cells.Item(25, 14).Value = 1.033555208843966
